$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace drive letter D: with F: in the path column (column A), rows 2-9
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = $current.Replace("D:/", "F:/")
    }
}

# Move selection to B15 as in the edited file
$ws.Range("B15").Select()

$wb.Save()
